$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for row 2 moved forward ~1 minute.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 05:02:16"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for row 2.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 05:02:11"
$wsZhCn.Range("K2").Value = "2016-08-26 05:02:32"

# de-de sheet: "Correspond Handback DateTime" for row 2.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-26 05:02:39"
